$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the formatting of the existing
# header cells (bold, centered, bordered - style used by B1:H1).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I ("I0") and J ("IF"), rows 2-46.
$data = @(
    @{Row=2; I=7; J=8},
    @{Row=3; I=7; J=7},
    @{Row=4; I=8; J=9},
    @{Row=5; I=7; J=7},
    @{Row=6; I=8; J=8},
    @{Row=7; I=6; J=7},
    @{Row=8; I=7; J=7},
    @{Row=9; I=6; J=7},
    @{Row=10; I=9; J=9},
    @{Row=11; I=6; J=7},
    @{Row=12; I=4; J=5},
    @{Row=13; I=4; J=5},
    @{Row=14; I=7; J=7},
    @{Row=15; I=6; J=6},
    @{Row=16; I=5; J=6},
    @{Row=17; I=1; J=2},
    @{Row=18; I=5; J=6},
    @{Row=19; I=5; J=6},
    @{Row=20; I=4; J=6},
    @{Row=21; I=6; J=6},
    @{Row=22; I=5; J=5},
    @{Row=23; I=6; J=6},
    @{Row=24; I=4; J=5},
    @{Row=25; I=6; J=6},
    @{Row=26; I=6; J=7},
    @{Row=27; I=6; J=7},
    @{Row=28; I=1; J=2},
    @{Row=29; I=7; J=8},
    @{Row=30; I=3; J=4},
    @{Row=31; I=7; J=9},
    @{Row=32; I=2; J=3},
    @{Row=33; I=7; J=7},
    @{Row=34; I=5; J=6},
    @{Row=35; I=8; J=9},
    @{Row=36; I=4; J=5},
    @{Row=37; I=5; J=6},
    @{Row=38; I=8; J=8},
    @{Row=39; I=7; J=7},
    @{Row=40; I=8; J=8},
    @{Row=41; I=5; J=6},
    @{Row=42; I=8; J=8},
    @{Row=43; I=5; J=5},
    @{Row=44; I=3; J=4},
    @{Row=45; I=9; J=9},
    @{Row=46; I=7; J=7}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value = $entry.I
    $ws.Cells.Item($r, 10).Value = $entry.J
}
